$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "durationfinal"
$ws.Range("C18").Value = "14.23 [5.82,23.53]"
$ws.Range("D18").Value = "0 (0%)"
$ws.Range("E18").Value = "14.4 [13.66,18.1]"
$ws.Range("F18").Value = "0 (0%)"
